$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-22 20:18:24"
$ws.Range("O2").Value = "6.2 °C"
$ws.Range("E3").Value = "2026-02-22 20:18:26"
$ws.Range("E4").Value = "2026-02-22 20:18:29"
$ws.Range("O4").Value = "12.5 °C"
$ws.Range("E5").Value = "2026-02-22 20:18:31"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "28%"
$ws.Range("H5").ClearFormats()
$ws.Range("H5").Borders.LineStyle = 1
$ws.Range("E6").Value = "2026-02-22 20:18:34"
$ws.Range("O6").Value = "13.1 °C"
$ws.Range("E7").Value = "2026-02-22 20:18:36"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "59%"
$ws.Range("H7").ClearFormats()
$ws.Range("H7").Borders.LineStyle = 1
$ws.Range("E8").Value = "2026-02-22 20:18:39"
$ws.Range("E9").Value = "2026-02-22 20:18:41"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "75%"
$ws.Range("H9").ClearFormats()
$ws.Range("H9").Borders.LineStyle = 1
$ws.Range("E10").Value = "2026-02-22 20:18:44"
$ws.Range("E11").Value = "2026-02-22 20:18:46"
$ws.Range("E12").Value = "2026-02-22 20:18:48"
$ws.Range("E13").Value = "2026-02-22 20:18:51"
$ws.Range("J13").Value = "1030.2 hPa"
$ws.Range("E14").Value = "2026-02-22 20:18:53"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "72%"
$ws.Range("H14").ClearFormats()
$ws.Range("H14").Borders.LineStyle = 1
$ws.Range("O14").Value = "12.1 °C"
$ws.Range("E15").Value = "2026-02-22 20:18:56"
$ws.Range("E16").Value = "2026-02-22 20:18:58"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "18%"
$ws.Range("H16").ClearFormats()
$ws.Range("H16").Borders.LineStyle = 1
$ws.Range("K16").Value = "13.0 MJ/m2"
$ws.Range("E17").Value = "2026-02-22 20:19:00"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "29%"
$ws.Range("H17").ClearFormats()
$ws.Range("H17").Borders.LineStyle = 1
$ws.Range("O17").Value = "10.1 °C"
$ws.Range("E18").Value = "2026-02-22 20:19:03"
$ws.Range("O18").Value = "10.2 °C"
$ws.Range("E19").Value = "2026-02-22 20:19:05"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "46%"
$ws.Range("H19").ClearFormats()
$ws.Range("H19").Borders.LineStyle = 1
$ws.Range("E20").Value = "2026-02-22 20:19:08"
$ws.Range("E21").Value = "2026-02-22 20:19:10"
$ws.Range("E22").Value = "2026-02-22 20:19:13"
$ws.Range("E23").Value = "2026-02-22 20:19:15"
$ws.Range("E24").Value = "2026-02-22 20:19:17"
$ws.Range("J24").Value = "1029.6 hPa"
$ws.Range("E25").Value = "2026-02-22 20:19:20"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "27%"
$ws.Range("H25").ClearFormats()
$ws.Range("H25").Borders.LineStyle = 1
$ws.Range("E26").Value = "2026-02-22 20:19:22"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "36%"
$ws.Range("H26").ClearFormats()
$ws.Range("H26").Borders.LineStyle = 1
$ws.Range("J26").Value = "1026.0 hPa"
$ws.Range("E27").Value = "2026-02-22 20:19:25"
$ws.Range("O27").Value = "6.7 °C"
$ws.Range("E28").Value = "2026-02-22 20:19:27"
$ws.Range("E29").Value = "2026-02-22 20:19:30"
$ws.Range("E30").Value = "2026-02-22 20:19:32"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "72%"
$ws.Range("H30").ClearFormats()
$ws.Range("H30").Borders.LineStyle = 1
$ws.Range("O30").Value = "12.3 °C"
$ws.Range("E31").Value = "2026-02-22 20:19:35"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "62%"
$ws.Range("H31").ClearFormats()
$ws.Range("H31").Borders.LineStyle = 1
$ws.Range("J31").Value = "1026.5 hPa"
$ws.Range("E32").Value = "2026-02-22 20:19:37"
$ws.Range("O32").Value = "6.2 °C"
$ws.Range("E33").Value = "2026-02-22 20:19:40"
$ws.Range("E34").Value = "2026-02-22 20:19:42"
$ws.Range("E35").Value = "2026-02-22 20:19:45"
$ws.Range("L35").Value = "32.0 km/h - 287º 19:52 TU"
$ws.Range("E36").Value = "2026-02-22 20:19:47"
$ws.Range("K36").Value = "15.1 MJ/m2"
$ws.Range("E37").Value = "2026-02-22 20:19:50"
$ws.Range("E38").Value = "2026-02-22 20:19:52"
$ws.Range("E39").Value = "2026-02-22 20:19:55"
$ws.Range("L39").Value = "25.6 km/h - 308º 19:53 TU"
$ws.Range("O39").Value = "5.2 °C"
$ws.Range("E40").Value = "2026-02-22 20:19:57"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "55%"
$ws.Range("H40").ClearFormats()
$ws.Range("H40").Borders.LineStyle = 1
$ws.Range("O40").Value = "10.2 °C"
$ws.Range("E41").Value = "2026-02-22 20:20:00"
$ws.Range("O41").Value = "11.3 °C"
$ws.Range("E42").Value = "2026-02-22 20:20:02"
$ws.Range("E43").Value = "2026-02-22 20:20:04"
$ws.Range("E44").Value = "2026-02-22 20:20:07"
$ws.Range("E45").Value = "2026-02-22 20:20:09"
$ws.Range("E46").Value = "2026-02-22 20:20:11"
